# fix(testcase): fix testcase generator some issue
# The "模块名" template sheet has a merged header cell (K1:O1) that labelled
# the second block of columns with the generic text "测试结果". The
# generator is being fixed to stamp the round number into that header, so
# the template default is updated to read "第六轮测试结果" instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("模块名")

# Rename the merged header cell K1:O1 from the generic "测试结果" label to
# "第六轮测试结果", matching the K1:O1 row-2 sibling headers' formatting
# (plain header style, no special font override).
$k1 = $ws.Range("K1")
$k1.Value = "第六轮测试结果"

# Match the normal header styling used by the rest of row 1 (e.g. L1),
# clearing the stray font override that the old duplicate style carried.
$l1 = $ws.Range("L1")
$k1.Font.Name = $l1.Font.Name
$k1.Font.Size = $l1.Font.Size
$k1.Font.Bold = $l1.Font.Bold
$k1.Font.Italic = $l1.Font.Italic

# Move the active selection to O3, matching where the author's cursor
# ended up after making the edit.
$ws.Range("O3").Select()
